# Update the cryptos list with the latest scraped prices/volumes
# (mirrors the GitHub Actions scraper commit on Thu Nov 23 05:40:40 UTC 2023)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numbers as plain text (e.g. "37.406.80" uses
# dots as thousands separators, and values like "14.60" / "0.780" rely on
# trailing zeros). Force those particular cells to Text format first so
# Excel doesn't silently reinterpret them as numbers and strip the
# formatting-significant digits.
$textCells = @(
    "D5","D7","D9","D10","D11","D14","D15","D16","D17","D20","D21","D23",
    "D27","D28","D30","D33","D39","D40","D41","D44","D45","D47","D48"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "37.406.80"
$ws.Range("E2").Value = "  +2.80%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.066.56"
$ws.Range("E3").Value = "  +4.05%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.05%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "235.43"
$ws.Range("E5").Value = "  -0.23%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  +2.88%  "

# --- Row 7: Solana ---
$ws.Range("D7").Value = "58.29"
$ws.Range("E7").Value = "  +7.07%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  -0.03%  "

# --- Row 9: Cardano ---
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +3.53%  "

# --- Row 10: OKB ---
$ws.Range("D10").Value = "58.82"
$ws.Range("E10").Value = "  +1.43%  "

# --- Row 11: Dogecoin ---
$ws.Range("D11").Value = "0.0762"
$ws.Range("E11").Value = "  +1.82%  "

# --- Row 12: TRON ---
$ws.Range("E12").Value = "  +2.94%  "

# --- Row 13: WrappedliquidstakedEther2.0 ---
$ws.Range("D13").Value = "2.371.50"
$ws.Range("E13").Value = "  +4.01%  "

# --- Row 14: Chainlink ---
$ws.Range("D14").Value = "14.60"
$ws.Range("E14").Value = "  +3.16%  "

# --- Row 15: Avalanche ---
$ws.Range("D15").Value = "20.98"
$ws.Range("E15").Value = "  +4.03%  "

# --- Row 16: Polygon ---
$ws.Range("D16").Value = "0.780"
$ws.Range("E16").Value = "  +3.33%  "

# --- Row 17: Polkadot ---
$ws.Range("D17").Value = "5.18"
$ws.Range("E17").Value = "  +2.70%  "

# --- Row 18: WrappedEther ---
$ws.Range("D18").Value = "2.055.87"
$ws.Range("E18").Value = "  +3.43%  "

# --- Row 19: WrappedBTC ---
$ws.Range("D19").Value = "37.603.39"
$ws.Range("E19").Value = "  +3.41%  "

# --- Row 20: Uniswap ---
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  +17.53%  "

# --- Row 21: Litecoin ---
$ws.Range("D21").Value = "69.06"
$ws.Range("E21").Value = "  +1.88%  "

# --- Row 22: ShibaInu ---
$ws.Range("D22").Value = "0.0₃0816"
$ws.Range("E22").Value = "  +1.65%  "

# --- Row 23: BitcoinCash ---
$ws.Range("D23").Value = "225.98"
$ws.Range("E23").Value = "  +2.14%  "

# --- Row 25: PancakeSwap ---
$ws.Range("E25").Value = "  +2.40%  "

# --- Row 26: Toncoin ---
$ws.Range("E26").Value = "  +1.02%  "

# --- Row 27 & 28: ImmutableX / Monero swap order (Monero now ranks 27, ImmutableX 28) ---
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "164.66"
$ws.Range("E27").Value = "  +1.10%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "1.51"
$ws.Range("E28").Value = "  +14.19%  "

# --- Row 29: Cosmos ---
$ws.Range("E29").Value = "  +2.67%  "

# --- Row 30: EthereumClassic ---
$ws.Range("D30").Value = "19.18"

# --- Row 31: Kaspa ---
$ws.Range("E31").Value = "  -0.38%  "

# --- Row 32: Stellar ---
$ws.Range("E32").Value = "  +2.08%  "

# --- Row 33: Filecoin ---
$ws.Range("D33").Value = "4.50"
$ws.Range("E33").Value = "  +2.85%  "

# --- Row 34: Hedera ---
$ws.Range("E34").Value = "  +2.97%  "

# --- Row 35: LidoDAOToken ---
$ws.Range("E35").Value = "  +9.34%  "

# --- Row 36: InternetComputer(DFINITY) ---
$ws.Range("E36").Value = "  +6.41%  "

# --- Row 37: RenderToken ---
$ws.Range("E37").Value = "  +3.32%  "

# --- Row 38: BinanceUSD ---
$ws.Range("E38").Value = "  +0.04%  "

# --- Row 39: WEMIXToken ---
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  +0.77%  "

# --- Row 40: THORChain ---
$ws.Range("D40").Value = "5.88"
$ws.Range("E40").Value = "  +7.47%  "

# --- Row 41: Cronos ---
$ws.Range("D41").Value = "0.0986"
$ws.Range("E41").Value = "  +7.39%  "

# --- Row 42: HuobiToken ---
$ws.Range("E42").Value = "  -0.93%  "

# --- Rows 43-45: FTXToken / Maker / Aave re-ranked ---
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.476.10"
$ws.Range("E43").Value = "  +1.77%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "97.04"
$ws.Range("E44").Value = "  +8.52%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "4.35"
$ws.Range("E45").Value = "  +20.18%  "

# --- Row 46: TrustWalletToken ---
$ws.Range("E46").Value = "  +6.50%  "

# --- Row 47: VeChain ---
$ws.Range("D47").Value = "0.0210"
$ws.Range("E47").Value = "  +4.36%  "

# --- Row 48: InjectiveProtocol ---
$ws.Range("D48").Value = "15.99"
$ws.Range("E48").Value = "  +6.91%  "

# --- Row 49: ARBITRUM ---
$ws.Range("E49").Value = "  +4.23%  "

# --- Row 50: FraxShare ---
$ws.Range("E50").Value = "  +6.89%  "

# --- Row 51: MXToken ---
$ws.Range("E51").Value = "  +2.22%  "
